# Applies scheduled-runner Leve-profit cell updates across all job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) in the Golem_Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 44164.25
$ws.Range("J3").Value = 44164.25
$ws.Range("L3").Value = 44164.25
$ws.Range("N3").Value = -44392.25
$ws.Range("H32").Value = 6500.5
$ws.Range("J32").Value = 6500.5
$ws.Range("L32").Value = 6500.5
$ws.Range("N32").Value = -7152.5
$ws.Range("H55").Value = 1878
$ws.Range("I55").Value = 3699.5
$ws.Range("J55").Value = 663.6667
$ws.Range("K55").Value = 3699.5
$ws.Range("L55").Value = 663.6667
$ws.Range("M55").Value = -3485.5
$ws.Range("N55").Value = -1091.6667
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = $null
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = $null
$ws.Range("H88").Value = 5000
$ws.Range("J88").Value = 5000
$ws.Range("L88").Value = 5000
$ws.Range("N88").Value = -5812
$ws.Range("H91").Value = 5000
$ws.Range("J91").Value = 5000
$ws.Range("L91").Value = 5000
$ws.Range("N91").Value = -7808
$ws.Range("H100").Value = 1398.8334
$ws.Range("I100").Value = 1398.8334
$ws.Range("K100").Value = 1398.8334
$ws.Range("M100").Value = -857.8334
$ws.Range("H102").Value = 44164.25
$ws.Range("J102").Value = 44164.25
$ws.Range("L102").Value = 44164.25
$ws.Range("N102").Value = -50654.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1636.3334
$ws.Range("I97").Value = 1636.3334
$ws.Range("K97").Value = 1636.3334
$ws.Range("M97").Value = -1140.3334
$ws.Range("H132").Value = 2377.5
$ws.Range("I132").Value = 1004
$ws.Range("K132").Value = 3012
$ws.Range("M132").Value = -482

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 1316.6666
$ws.Range("I29").Value = 1525
$ws.Range("J29").Value = 900
$ws.Range("K29").Value = 1525
$ws.Range("L29").Value = 900
$ws.Range("M29").Value = -1236
$ws.Range("N29").Value = -1478
$ws.Range("H64").Value = 462.6
$ws.Range("I64").Value = 375.42856
$ws.Range("K64").Value = 375.42856
$ws.Range("M64").Value = -150.42856
$ws.Range("H67").Value = 462.6
$ws.Range("I67").Value = 375.42856
$ws.Range("K67").Value = 375.42856
$ws.Range("M67").Value = 404.57144
$ws.Range("H86").Value = 2040.2
$ws.Range("I86").Value = 1800.25
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 1800.25
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -677.25
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 2040.2
$ws.Range("I89").Value = 1800.25
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 9001.25
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -3385.25
$ws.Range("N89").Value = -26232
$ws.Range("H94").Value = 2407.6924
$ws.Range("I94").Value = 1885.7142
$ws.Range("K94").Value = 1885.7142
$ws.Range("M94").Value = -1434.7142
$ws.Range("H99").Value = 4745
$ws.Range("I99").Value = 5024.875
$ws.Range("K99").Value = 5024.875
$ws.Range("M99").Value = -3526.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 984.1429000000001
$ws.Range("J22").Value = 1399.75
$ws.Range("L22").Value = 1399.75
$ws.Range("N22").Value = -2099.75
$ws.Range("H31").Value = 4992.4707
$ws.Range("I31").Value = 4324.8
$ws.Range("K31").Value = 4324.8
$ws.Range("M31").Value = -4029.8
$ws.Range("H34").Value = 4992.4707
$ws.Range("I34").Value = 4324.8
$ws.Range("K34").Value = 4324.8
$ws.Range("M34").Value = -4122.8
$ws.Range("H36").Value = 4800
$ws.Range("I36").Value = 4500
$ws.Range("J36").Value = 5100
$ws.Range("K36").Value = 4500
$ws.Range("L36").Value = 5100
$ws.Range("M36").Value = -4112
$ws.Range("N36").Value = -5876
$ws.Range("H40").Value = 4800
$ws.Range("I40").Value = 4500
$ws.Range("J40").Value = 5100
$ws.Range("K40").Value = 4500
$ws.Range("L40").Value = 5100
$ws.Range("M40").Value = -4340
$ws.Range("N40").Value = -5420

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2592.75
$ws.Range("I68").Value = 859.5
$ws.Range("K68").Value = 2578.5
$ws.Range("M68").Value = -1767.5
$ws.Range("H71").Value = 2592.75
$ws.Range("I71").Value = 859.5
$ws.Range("K71").Value = 7735.5
$ws.Range("M71").Value = -3679.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 4292279.5
$ws.Range("I11").Value = 8862500
$ws.Range("J11").Value = 229861
$ws.Range("K11").Value = 8862500
$ws.Range("L11").Value = 229861
$ws.Range("M11").Value = -8862361
$ws.Range("N11").Value = -230139
$ws.Range("H12").Value = 1499
$ws.Range("J12").Value = 1499
$ws.Range("L12").Value = 1499
$ws.Range("N12").Value = -1779
$ws.Range("H19").Value = 5006
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").Value = $null
$ws.Range("H97").Value = 1630.4
$ws.Range("I97").Value = 1663
$ws.Range("K97").Value = 1663
$ws.Range("M97").Value = -1167

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 893.125
$ws.Range("I46").Value = 749.75
$ws.Range("J46").Value = 1036.5
$ws.Range("K46").Value = 749.75
$ws.Range("L46").Value = 1036.5
$ws.Range("M46").Value = -561.75
$ws.Range("N46").Value = -1412.5
$ws.Range("H61").Value = 1693
$ws.Range("J61").Value = 1001
$ws.Range("L61").Value = 1001
$ws.Range("N61").Value = -1405
$ws.Range("H100").Value = 5000
$ws.Range("I100").Value = 5000
$ws.Range("K100").Value = 5000
$ws.Range("M100").Value = -4459
$ws.Range("H113").Value = 1693
$ws.Range("J113").Value = 1001
$ws.Range("L113").Value = 1001
$ws.Range("N113").Value = -5341
$ws.Range("H132").Value = 1498
$ws.Range("I132").Value = 1498
$ws.Range("K132").Value = 4494
$ws.Range("M132").Value = -1964
$ws.Range("H135").Value = 99995
$ws.Range("J135").Value = 99995
$ws.Range("L135").Value = 99995
$ws.Range("N135").Value = -110135

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 42500
$ws.Range("J40").Value = 25000
$ws.Range("L40").Value = 25000
$ws.Range("N40").Value = -25298
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = $null
$ws.Range("H54").Value = 30000
$ws.Range("J54").Value = 30000
$ws.Range("L54").Value = 30000
$ws.Range("N54").Value = -31040
$ws.Range("H93").Value = 50000
$ws.Range("I93").Value = 50000
$ws.Range("K93").Value = 50000
$ws.Range("M93").Value = -47504
$ws.Range("H104").Value = 14374.25
$ws.Range("J104").Value = 14374.25
$ws.Range("L104").Value = 14374.25
$ws.Range("N104").Value = -21362.25
$ws.Range("H126").Value = 1166.2
$ws.Range("I126").Value = 1166.2
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3498.6
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1028.6
$ws.Range("N126").Value = $null

